# Scheduled-runner refresh of cached Kraken market data across the
# per-job Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ
# (H/I/J), LevePriceNQ / LevePriceHQ (K/L) and the derived profit columns
# (M/N) for the rows whose market snapshot changed; cells that no longer
# have a computed profit are cleared back out.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 148.66667  # H6
$ws.Cells.Item(6, 10).Value = 200  # J6
$ws.Cells.Item(6, 12).Value = 600  # L6
$ws.Cells.Item(6, 14).Value = -824  # N6
$ws.Cells.Item(28, 8).Value = 100000  # H28
$ws.Cells.Item(28, 9).Value = 100000  # I28
$ws.Cells.Item(28, 11).Value = 100000  # K28
$ws.Cells.Item(28, 13).Value = -99515  # M28
$ws.Cells.Item(33, 8).Value = 410.66666  # H33
$ws.Cells.Item(33, 9).Value = 87.53846  # I33
$ws.Cells.Item(33, 10).Value = 2511  # J33
$ws.Cells.Item(33, 11).Value = 87.53846  # K33
$ws.Cells.Item(33, 12).Value = 2511  # L33
$ws.Cells.Item(33, 13).Value = 141.46154  # M33
$ws.Cells.Item(33, 14).Value = -2969  # N33
$ws.Cells.Item(51, 8).Value = 2400  # H51
$ws.Cells.Item(51, 9).Value = 2400  # I51
$ws.Cells.Item(51, 11).Value = 2400  # K51
$ws.Cells.Item(51, 13).Value = -1916  # M51
$ws.Cells.Item(68, 8).Value = 0  # H68
$ws.Cells.Item(68, 9).Value = 0  # I68
$ws.Cells.Item(68, 11).Value = 0  # K68
$ws.Cells.Item(68, 13).Value = $null  # M68
$ws.Cells.Item(71, 8).Value = 0  # H71
$ws.Cells.Item(71, 9).Value = 0  # I71
$ws.Cells.Item(71, 11).Value = 0  # K71
$ws.Cells.Item(71, 13).Value = $null  # M71
$ws.Cells.Item(106, 8).Value = 2856.4285  # H106
$ws.Cells.Item(106, 9).Value = 1995  # I106
$ws.Cells.Item(106, 11).Value = 1995  # K106
$ws.Cells.Item(106, 13).Value = -1364  # M106
$ws.Cells.Item(113, 8).Value = 2000  # H113
$ws.Cells.Item(113, 9).Value = 2000  # I113
$ws.Cells.Item(113, 11).Value = 2000  # K113
$ws.Cells.Item(113, 13).Value = 1254  # M113
$ws.Cells.Item(132, 8).Value = 2617.4075  # H132
$ws.Cells.Item(132, 9).Value = 2333.6924  # I132
$ws.Cells.Item(132, 11).Value = 7001.0772  # K132
$ws.Cells.Item(132, 13).Value = -4471.0772  # M132
$ws.Cells.Item(137, 8).Value = 2416.6667  # H137
$ws.Cells.Item(137, 9).Value = 2416.6667  # I137
$ws.Cells.Item(137, 11).Value = 7250.000100000001  # K137
$ws.Cells.Item(137, 13).Value = -4700.000100000001  # M137
$ws.Cells.Item(138, 8).Value = 3239.6667  # H138
$ws.Cells.Item(138, 10).Value = 3943.6667  # J138
$ws.Cells.Item(138, 12).Value = 11831.0001  # L138
$ws.Cells.Item(138, 14).Value = -22111.0001  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 613  # H2
$ws.Cells.Item(2, 9).Value = 0  # I2
$ws.Cells.Item(2, 10).Value = 613  # J2
$ws.Cells.Item(2, 11).Value = 0  # K2
$ws.Cells.Item(2, 12).Value = $null  # L2
$ws.Cells.Item(2, 13).Value = 613  # M2
$ws.Cells.Item(2, 14).Value = -839  # N2
$ws.Cells.Item(32, 8).Value = 3634.875  # H32
$ws.Cells.Item(32, 9).Value = 4180  # I32
$ws.Cells.Item(32, 11).Value = 4180  # K32
$ws.Cells.Item(32, 13).Value = -3893  # M32
$ws.Cells.Item(116, 8).Value = 613  # H116
$ws.Cells.Item(116, 9).Value = 0  # I116
$ws.Cells.Item(116, 10).Value = 613  # J116
$ws.Cells.Item(116, 11).Value = 0  # K116
$ws.Cells.Item(116, 12).Value = $null  # L116
$ws.Cells.Item(116, 13).Value = 613  # M116
$ws.Cells.Item(116, 14).Value = -5201  # N116
$ws.Cells.Item(137, 8).Value = 74997  # H137
$ws.Cells.Item(137, 9).Value = 49999  # I137
$ws.Cells.Item(137, 10).Value = 99995  # J137
$ws.Cells.Item(137, 11).Value = 49999  # K137
$ws.Cells.Item(137, 12).Value = 99995  # L137
$ws.Cells.Item(137, 13).Value = -44899  # M137
$ws.Cells.Item(137, 14).Value = -110195  # N137

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 613  # H3
$ws.Cells.Item(3, 9).Value = 0  # I3
$ws.Cells.Item(3, 10).Value = 613  # J3
$ws.Cells.Item(3, 11).Value = 0  # K3
$ws.Cells.Item(3, 12).Value = $null  # L3
$ws.Cells.Item(3, 13).Value = 613  # M3
$ws.Cells.Item(3, 14).Value = -841  # N3
$ws.Cells.Item(64, 8).Value = 796  # H64
$ws.Cells.Item(64, 9).Value = 255.2  # I64
$ws.Cells.Item(64, 10).Value = 3500  # J64
$ws.Cells.Item(64, 11).Value = 255.2  # K64
$ws.Cells.Item(64, 12).Value = 3500  # L64
$ws.Cells.Item(64, 13).Value = -30.19999999999999  # M64
$ws.Cells.Item(64, 14).Value = -3950  # N64
$ws.Cells.Item(67, 8).Value = 796  # H67
$ws.Cells.Item(67, 9).Value = 255.2  # I67
$ws.Cells.Item(67, 10).Value = 3500  # J67
$ws.Cells.Item(67, 11).Value = 255.2  # K67
$ws.Cells.Item(67, 12).Value = 3500  # L67
$ws.Cells.Item(67, 13).Value = 524.8  # M67
$ws.Cells.Item(67, 14).Value = -5060  # N67
$ws.Cells.Item(76, 8).Value = 17499  # H76
$ws.Cells.Item(76, 10).Value = 17499  # J76
$ws.Cells.Item(76, 12).Value = 17499  # L76
$ws.Cells.Item(76, 14).Value = -18129  # N76
$ws.Cells.Item(79, 8).Value = 17499  # H79
$ws.Cells.Item(79, 10).Value = 17499  # J79
$ws.Cells.Item(79, 12).Value = 17499  # L79
$ws.Cells.Item(79, 14).Value = -19683  # N79
$ws.Cells.Item(107, 9).Value = 666.6667  # I107
$ws.Cells.Item(107, 11).Value = 666.6667  # K107
$ws.Cells.Item(107, 13).Value = 1253.3333  # M107
$ws.Cells.Item(125, 8).Value = 99995  # H125
$ws.Cells.Item(125, 10).Value = 99995  # J125
$ws.Cells.Item(125, 12).Value = 99995  # L125
$ws.Cells.Item(125, 14).Value = -109835  # N125

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(35, 8).Value = 663.2  # H35
$ws.Cells.Item(35, 9).Value = 824.5  # I35
$ws.Cells.Item(35, 10).Value = 18  # J35
$ws.Cells.Item(35, 11).Value = 824.5  # K35
$ws.Cells.Item(35, 12).Value = 18  # L35
$ws.Cells.Item(35, 13).Value = -530.5  # M35
$ws.Cells.Item(35, 14).Value = -606  # N35

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 429043.22  # H4
$ws.Cells.Item(4, 9).Value = 417175.5  # I4
$ws.Cells.Item(4, 11).Value = 1251526.5  # K4
$ws.Cells.Item(4, 13).Value = -1251414.5  # M4
$ws.Cells.Item(34, 8).Value = 1116.7778  # H34
$ws.Cells.Item(34, 9).Value = 683.6667  # I34
$ws.Cells.Item(34, 11).Value = 2051.0001  # K34
$ws.Cells.Item(34, 13).Value = -1967.0001  # M34
$ws.Cells.Item(114, 8).Value = 757.75  # H114
$ws.Cells.Item(114, 9).Value = 0  # I114
$ws.Cells.Item(114, 10).Value = 757.75  # J114
$ws.Cells.Item(114, 11).Value = 0  # K114
$ws.Cells.Item(114, 12).Value = $null  # L114
$ws.Cells.Item(114, 13).Value = 2273.25  # M114
$ws.Cells.Item(114, 14).Value = -8781.25  # N114
$ws.Cells.Item(122, 8).Value = 494  # H122
$ws.Cells.Item(122, 10).Value = 0  # J122
$ws.Cells.Item(122, 12).Value = 0  # L122
$ws.Cells.Item(122, 14).Value = $null  # N122
$ws.Cells.Item(123, 8).Value = 2000  # H123
$ws.Cells.Item(123, 9).Value = 2000  # I123
$ws.Cells.Item(123, 11).Value = 6000  # K123
$ws.Cells.Item(123, 13).Value = -3550  # M123

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 9).Value = 19000  # I5
$ws.Cells.Item(5, 11).Value = 19000  # K5
$ws.Cells.Item(5, 13).Value = -18888  # M5
$ws.Cells.Item(39, 8).Value = 0  # H39
$ws.Cells.Item(39, 10).Value = 0  # J39
$ws.Cells.Item(39, 12).Value = $null  # L39
$ws.Cells.Item(39, 14).Value = 0  # N39
$ws.Cells.Item(132, 8).Value = 3174.375  # H132
$ws.Cells.Item(132, 9).Value = 2649.1667  # I132
$ws.Cells.Item(132, 11).Value = 7947.500100000001  # K132
$ws.Cells.Item(132, 13).Value = -5417.500100000001  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 8344168  # H3
$ws.Cells.Item(3, 10).Value = 12500  # J3
$ws.Cells.Item(3, 12).Value = 12500  # L3
$ws.Cells.Item(3, 14).Value = -12724  # N3
$ws.Cells.Item(11, 8).Value = 0  # H11
$ws.Cells.Item(11, 9).Value = 0  # I11
$ws.Cells.Item(11, 11).Value = 0  # K11
$ws.Cells.Item(11, 13).Value = $null  # M11
$ws.Cells.Item(15, 8).Value = 8344168  # H15
$ws.Cells.Item(15, 10).Value = 12500  # J15
$ws.Cells.Item(15, 12).Value = 12500  # L15
$ws.Cells.Item(15, 14).Value = -12840  # N15
$ws.Cells.Item(25, 8).Value = 5000  # H25
$ws.Cells.Item(25, 9).Value = 5000  # I25
$ws.Cells.Item(25, 11).Value = 5000  # K25
$ws.Cells.Item(25, 13).Value = -4770  # M25
$ws.Cells.Item(93, 8).Value = 1875.875  # H93
$ws.Cells.Item(93, 9).Value = 1826.25  # I93
$ws.Cells.Item(93, 11).Value = 1826.25  # K93
$ws.Cells.Item(93, 13).Value = -578.25  # M93
$ws.Cells.Item(122, 8).Value = 4259.2  # H122
$ws.Cells.Item(122, 9).Value = 4259.2  # I122
$ws.Cells.Item(122, 11).Value = 12777.6  # K122
$ws.Cells.Item(122, 13).Value = -10327.6  # M122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(19, 8).Value = 0  # H19
$ws.Cells.Item(19, 9).Value = 0  # I19
$ws.Cells.Item(19, 11).Value = 0  # K19
$ws.Cells.Item(19, 13).Value = $null  # M19
$ws.Cells.Item(37, 8).Value = 30026  # H37
$ws.Cells.Item(37, 9).Value = 30026  # I37
$ws.Cells.Item(37, 11).Value = 30026  # K37
$ws.Cells.Item(37, 13).Value = -29823  # M37
$ws.Cells.Item(62, 8).Value = 4500  # H62
$ws.Cells.Item(62, 9).Value = 4500  # I62
$ws.Cells.Item(62, 11).Value = 4500  # K62
$ws.Cells.Item(62, 13).Value = -3876  # M62
$ws.Cells.Item(65, 8).Value = 4500  # H65
$ws.Cells.Item(65, 9).Value = 4500  # I65
$ws.Cells.Item(65, 11).Value = 22500  # K65
$ws.Cells.Item(65, 13).Value = -19380  # M65
$ws.Cells.Item(107, 8).Value = 2625.25  # H107
$ws.Cells.Item(107, 9).Value = 1600.4  # I107
$ws.Cells.Item(107, 11).Value = 4801.200000000001  # K107
$ws.Cells.Item(107, 13).Value = -2881.200000000001  # M107
